$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C152").Value = 1.238269994140582
$ws.Range("C153").Value = 1.226948046732958
$ws.Range("C154").Value = 1.226648974683219
$ws.Range("C155").Value = 1.251598820715003
$ws.Range("C156").Value = 1.239783964768033
$ws.Range("C157").Value = 1.231341914339102
$ws.Range("C158").Value = 1.228113977519594
$ws.Range("C159").Value = 1.224596937211389
$ws.Range("C160").Value = 1.22466650227657
$ws.Range("C161").Value = 1.235753376957996
$ws.Range("C162").Value = 1.207394298574729
$ws.Range("C163").Value = 1.229196282763044
$ws.Range("C164").Value = 1.221464327955706
$ws.Range("C165").Value = 1.220614369390921
$ws.Range("C166").Value = 1.232256445635086
$ws.Range("C167").Value = 1.254837787206313
$ws.Range("C168").Value = 1.255368858611194
$ws.Range("C169").Value = 1.25329796961472
$ws.Range("C170").Value = 1.248173324851996
$ws.Range("C171").Value = 1.244793804339156
$ws.Range("C172").Value = 1.241089114964323
$ws.Range("C173").Value = 1.252630871270916
$ws.Range("C174").Value = 1.261138244488851
$ws.Range("C175").Value = 1.251640714470686
$ws.Range("C176").Value = 1.240701803075182
$ws.Range("C177").Value = 1.239882725373991
$ws.Range("C178").Value = 1.239229319877219
$ws.Range("C179").Value = 1.239385813920769
$ws.Range("C180").Value = 1.251827520035059
$ws.Range("C181").Value = 1.231903864561693
$ws.Range("C182").Value = 1.232453600326134
$ws.Range("C183").Value = 1.231111004911923
$ws.Range("C184").Value = 1.230371637033066
$ws.Range("C185").Value = 1.229065773288799
$ws.Range("C186").Value = 1.228182306305312
$ws.Range("C187").Value = 1.214760197372009
$ws.Range("C188").Value = 1.214737649731441
$ws.Range("C189").Value = 1.216393863917958
$ws.Range("C190").Value = 1.195224373986969
$ws.Range("C191").Value = 1.206114938435964
$ws.Range("C192").Value = 1.215058924347624
$ws.Range("C193").Value = 1.213940947482899
$ws.Range("C194").Value = 1.21513958041862
$ws.Range("C195").Value = 1.192994489511771
$ws.Range("C196").Value = 1.192555829941442
$ws.Range("C197").Value = 1.191883890304369
$ws.Range("C198").Value = 1.193851112784309
$ws.Range("C199").Value = 1.182983355163583
$ws.Range("C200").Value = 1.182850136424509
$ws.Range("C201").Value = 1.134412140823284
$ws.Range("C202").Value = 1.155892290501299
$ws.Range("C203").Value = 1.155293426993963
$ws.Range("C204").Value = 1.158447683640106
$ws.Range("C205").Value = 1.158721805390498
$ws.Range("C206").Value = 1.160352133350743
$ws.Range("C207").Value = 1.158309344632277
$ws.Range("C208").Value = 1.157563741048335
$ws.Range("C209").Value = 1.198220063701776
$ws.Range("C210").Value = 1.195588336915779
$ws.Range("C211").Value = 1.19465070692743
$ws.Range("C212").Value = 1.193187757764555
$ws.Range("C213").Value = 1.189550956518071
$ws.Range("C214").Value = 1.202302842221616
$ws.Range("C215").Value = 1.212014416051553
$ws.Range("C216").Value = 1.234030229593855
$ws.Range("C217").Value = 1.246098319938603
$ws.Range("C218").Value = 1.246273242820492
$ws.Range("C219").Value = 1.246651568432251
$ws.Range("C220").Value = 1.247129108394523
$ws.Range("C221").Value = 1.246475114943493
$ws.Range("C222").Value = 1.245643994462828
$ws.Range("C223").Value = 1.261691503158172
$ws.Range("C224").Value = 1.257156149604506
$ws.Range("C225").Value = 1.209662791130737
$ws.Range("C226").Value = 1.209156750441042
$ws.Range("C227").Value = 1.208776236756015
$ws.Range("C228").Value = 1.217877687926723
$ws.Range("C229").Value = 1.210319165117816
$ws.Range("C230").Value = 1.220889966517364
$ws.Range("C231").Value = 1.209350538239505
$ws.Range("C232").Value = 1.208327948655497
$ws.Range("C233").Value = 1.207513369285311
$ws.Range("C234").Value = 1.254963718808026
$ws.Range("C235").Value = 1.2541140304283
$ws.Range("C236").Value = 1.252689620250149
$ws.Range("C237").Value = 1.251961778778859
$ws.Range("C238").Value = 1.251876098844296
$ws.Range("C239").Value = 1.251671245669107
$ws.Range("C240").Value = 1.263847695985228
$ws.Range("C241").Value = 1.262663169287868
$ws.Range("C242").Value = 1.274977080183574
$ws.Range("C243").Value = 1.299960470038234
$ws.Range("C244").Value = 1.337909672269377
$ws.Range("C245").Value = 1.350402306315797
$ws.Range("C246").Value = 1.374317018280406
$ws.Range("C247").Value = 1.37166063009448
$ws.Range("C248").Value = 1.396536135417827
$ws.Range("C249").Value = 1.406209022047955
$ws.Range("C250").Value = 1.404521719344106
$ws.Range("C251").Value = 1.488161909805113
$ws.Range("C252").Value = 1.474467238416472
$ws.Range("C253").Value = 1.473574340200032
$ws.Range("C255").Value = 1.487336000930483
$ws.Range("C256").Value = 1.488512669176378
$ws.Range("C258").Value = 1.52837138557035
$ws.Range("C259").Value = 1.528552072797504
$ws.Range("C262").Value = 1.490000049940465
$ws.Range("C263").Value = 1.489791398203276
$ws.Range("C264").Value = 1.47429691590474
$ws.Range("C265").Value = 1.446347853925327
$ws.Range("C266").Value = 1.444053426552606
$ws.Range("C267").Value = 1.443412231066546
$ws.Range("C268").Value = 1.451868197636169
$ws.Range("C269").Value = 1.435880133875951
$ws.Range("C270").Value = 1.381077293843607
$ws.Range("C271").Value = 1.341679770778456
$ws.Range("C272").Value = 1.329227400069269
